# "Update máy mới về" - NhapXuatLaptop.xlsx
# Applies the business-data edits captured in the commit: a handful of
# sales rows that are now settled get hidden from the active view, a few
# machines that sold get their actual sale price + note filled in (some
# with a text note instead of a number, which turns the profit formula
# into #VALUE!), two "máy mới về" (newly arrived machines) rows get their
# extra cost filled in, and the running index column (A) is renumbered.
#
# NOTE: cell writes are done before the Hidden=$true toggles on purpose -
# writing into an already-hidden row causes the host to stamp a spurious
# autofit row height, which the source workbook does not have.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NhapXuat")
$ws.Activate()

# ---------------------------------------------------------------------
# Notes / buyers for rows that sold
# ---------------------------------------------------------------------
$ws.Range("K43").Value = "Chú Kiểm Nha Trang"
$ws.Range("N42").Value = "397531124140"
$ws.Range("K30").Value = "A Đoàn"
$ws.Range("K39").Value = "Hiếu"
$ws.Range("H42").Value = "6tr5"

# Row 33: sold for "7tr" (text note instead of numeric price -> formula errors)
$ws.Range("H33").Value = "7tr"

# Row 39: actual numeric sale price recorded + sale date
$ws.Range("H39").Value = 7000000
$ws.Range("P39").Value = 45265

# Row 41: sold for "7tr"
$ws.Range("H41").Value = "7tr"

# Row 42: "máy mới về" - extra cost recorded
$ws.Range("D42").Value = 850000
$ws.Range("D42").WrapText = $true

# Row 43: actual numeric sale price recorded
$ws.Range("H43").Value = 8500000

# Row 44: "máy mới về" - extra cost recorded + sold for "7tr"
$ws.Range("D44").Value = 850000
$ws.Range("H44").Value = "7tr"

# ---------------------------------------------------------------------
# Renumber the running index column now that row 39 is merged into the
# hidden/settled group
# ---------------------------------------------------------------------
$ws.Range("A39").WrapText = $true

$ws.Range("A40").Value = 39
$ws.Range("A40").WrapText = $true

$ws.Range("A41").Value = 40
$ws.Range("A41").WrapText = $true

$ws.Range("A42").Value = 41
$ws.Range("A42").WrapText = $true

$ws.Range("A43").Value = 42
$ws.Range("A43").WrapText = $true

$ws.Range("A44").Value = 43
$ws.Range("A44").WrapText = $true

# ---------------------------------------------------------------------
# Rows that are now sold / settled and should drop out of the active view
# (done last - see note above)
# ---------------------------------------------------------------------
22..27 | ForEach-Object { $ws.Rows.Item($_).Hidden = $true }
$ws.Rows.Item(30).Hidden = $true
$ws.Rows.Item(39).Hidden = $true
$ws.Rows.Item(43).Hidden = $true

# ---------------------------------------------------------------------
# Leave the selection on the cell the author was last looking at
# ---------------------------------------------------------------------
$ws.Range("H33").Select()
